$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry sets the cell to Text format first so that numeric-looking
# strings (e.g. "71.193.61", "1.00") are preserved exactly as text and
# are not auto-converted into numbers/dates by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.193.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.695.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +7.92%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.49"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.682.79"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +7.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.618"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +27.31%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "49.21"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.291.29"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "678.32"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.697.16"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.339.64"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.60"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.944"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.09"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.15"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.42%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.52"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.03"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "582.08"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.20"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.75"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0464"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +11.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.144"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.627.19"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.352"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0770"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.59"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.77"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.99%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.97"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.08%  "
